# Replace the Arabic-numeral verse-number cells in column A with their
# Hebrew-letter (gematria) equivalents, and right-align that column's
# verse-number cells (Excel introduces a new cell style with
# horizontal="right" for this).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = -4152
$xlBottom = -4107

$verseRows = @(2, 8, 21, 32, 42, 51, 58, 66, 73, 81)
$hebrewLetters = @("א", "ב", "ג", "ד", "ה", "ו", "ז", "ח", "ט", "י")

for ($i = 0; $i -lt $verseRows.Length; $i++) {
    $row = $verseRows[$i]
    $letter = $hebrewLetters[$i]

    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $letter
    $cell.HorizontalAlignment = $xlRight
    $cell.VerticalAlignment = $xlBottom
}
